# Generate Report for Handback
#
# The file "9f4ad892-dc37-4aec-b18d-4595c510be76.md" has been handed back
# (both zh-cn and de-de targets are now in sync with en-US). Update the
# status / handback-datetime columns on the "zh-cn" and "de-de" sheets for
# that file's row, and reflect the new status on the "Overview" sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: row for 9f4ad892-...md (row 3) -> zh-cn / de-de columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# zh-cn sheet: row for 9f4ad892-...md (row 3)
#   Status -> Handed back
#   Latest Handback DateTime -> new handback timestamp
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("G3").Value = "2016-02-22 09:23:23"

# ---------------------------------------------------------------------
# de-de sheet: row for 9f4ad892-...md (row 3)
#   Status -> Handed back
#   Latest Handback DateTime -> new handback timestamp
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("G3").Value = "2016-02-22 09:23:45"
